$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: status changes from "done" to "in Arbeit" with "Neutral" style
$ws.Range("B10").Value = "in Arbeit"
$ws.Range("B10").Style = "Neutral"

# New "Hinweis" column header
$ws.Range("D1").Value = "Hinweis"

# Make header row bold (A1:D1)
$ws.Range("A1:D1").Font.Bold = $true

# New note for the Cron-Job related todo row
$ws.Range("D6").Value = "Cron-Job"

# New row 11: addOrder umschreiben task
$ws.Range("A11").Value = "addOrder umschreiben"
$ws.Range("B11").Value = "in Arbeit"
$ws.Range("B11").Style = "Neutral"
$ws.Range("C11").Value = "Jonas"

$ws.Range("C11").Select() | Out-Null
